# edit.ps1 - applies the "Actualizado documento de requisitos" changes
# to the Acme Taxi "hackaton requirements" document.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $old"
    }
    return $ok
}

# ---------------------------------------------------------------
# 1. Drivers info bullet: drop "the maximum number of passengers
#    they can take at once, a number plate corresponding to their
#    vehicle, " from the list of stored driver data.
# ---------------------------------------------------------------
Replace-Text " a picture, the maximum number of passengers they can take at once, a number plate corresponding to their vehicle, their city of residence and their " " a picture, their city of residence and their "

# ---------------------------------------------------------------
# 2. Messaging paragraph: "can not" -> "cannot"
# ---------------------------------------------------------------
Replace-Text "These default folders can not be edited nor deleted" "These default folders cannot be edited nor deleted"

# ---------------------------------------------------------------
# 3. Messaging paragraph: clarify what blocking the sender means.
# ---------------------------------------------------------------
Replace-Text "which means the sender will no longer be able to send him a message" "which means the all future messages received written by that blocked sender will be automatically sent to the spam box"

# ---------------------------------------------------------------
# 5. Requests: an accepted request can be cancelled.
# ---------------------------------------------------------------
Replace-Text "it has not been taken by any driver." "it has not been taken by any driver. An accepted request can be cancelled as long as the moment when it was supposed to start has not passed"

# ---------------------------------------------------------------
# 6. Reviews: the related request must not be cancelled.
# ---------------------------------------------------------------
Replace-Text "this driver must have accepted a request from the user writing the review." "this driver must have accepted a request from the user writing the review. Said request must have not been cancelled and it must have already taken place"

# ---------------------------------------------------------------
# 7. Car registration wording: "Register car" -> "Register a car"
# ---------------------------------------------------------------
Replace-Text "Register car and associate" "Register a car and associate"

# ---------------------------------------------------------------
# 9. Notification alert wording simplified.
# ---------------------------------------------------------------
Replace-Text "A notification alert must appear whenever any actor has unread messages in his “notification box”." "A notification alert must appear whenever any actor has unread messages."

# ---------------------------------------------------------------
# 10. Sponsor advertisement placement changed.
# ---------------------------------------------------------------
Replace-Text "A random advertisement approved by a system admin from a sponsor must be displayed in one of the sides of the webpage." "A random advertisement approved by a system admin from a sponsor must be displayed at the bottom of the website."

# ---------------------------------------------------------------
# 4 (bookmark). Move the "_GoBack" bookmark from the end of the
# "A driver without a car can not accept requests." paragraph to
# just after "...the subject and the body. " (before "Every actor
# has the following system folders...").
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$bmRng = $d.Content
$bmRng.Find.Execute("the subject and the body. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRng.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

Write-Output "done"
